$wb = $excel.ActiveWorkbook
$origActiveSheetName = $wb.ActiveSheet.Name

# Add the new "Result1" sheet after the last existing sheet ("Result")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Result1"

# Header row (row 1) - only A1:E1 are populated (mirrors the source export, which
# leaves the data column F without its own header label)
$ws.Range("A1").Value = "Symbol"
$ws.Range("B1").Value = "LTP"
$ws.Range("C1").Value = "High52W"
$ws.Range("D1").Value = "DiffFrom52WHigh"
$ws.Range("E1").Value = "DiffFrom52WLow"

# Data rows 2-51: Symbol, LTP, High52W, DiffFrom52WHigh, Low52W, DiffFrom52WLow
$data = New-Object "object[,]" 50,6
$data[0,0] = "SIEMENS"
$data[0,1] = 3074.8
$data[0,2] = 6740.0
$data[0,3] = 45.62017804154303
$data[0,4] = 2450.0
$data[0,5] = 125.50204081632654
$data[1,0] = "RECLTD"
$data[1,1] = 357.5
$data[1,2] = 544.7
$data[1,3] = 65.63245823389022
$data[1,4] = 330.95
$data[1,5] = 108.02235987309261
$data[2,0] = "MAZDOCK"
$data[2,1] = 2490.0
$data[2,2] = 3775.0
$data[2,3] = 65.96026490066225
$data[2,4] = 1918.05
$data[2,5] = 129.81934777508408
$data[3,0] = "SOLARINDS"
$data[3,1] = 12264.0
$data[3,2] = 17820.0
$data[3,3] = 68.82154882154883
$data[3,4] = 8482.5
$data[3,5] = 144.58001768346597
$data[4,0] = "BAJAJHFL"
$data[4,1] = 94.38
$data[4,2] = 136.96
$data[4,3] = 68.91063084112149
$data[4,4] = 92.1
$data[4,5] = 102.47557003257329
$data[5,0] = "LODHA"
$data[5,1] = 1061.0
$data[5,2] = 1531.0
$data[5,3] = 69.30111038536904
$data[5,4] = 1035.15
$data[5,5] = 102.49722262474036
$data[6,0] = "ENRIN"
$data[6,1] = 2560.0
$data[6,2] = 3625.0
$data[6,3] = 70.62068965517241
$data[6,4] = 2508.8
$data[6,5] = 102.0408163265306
$data[7,0] = "JSWENERGY"
$data[7,1] = 482.8
$data[7,2] = 674.0
$data[7,3] = 71.6320474777448
$data[7,4] = 418.75
$data[7,5] = 115.29552238805971
$data[8,0] = "NAUKRI"
$data[8,1] = 1334.0
$data[8,2] = 1825.78
$data[8,3] = 73.06466277426634
$data[8,4] = 1157.0
$data[8,5] = 115.2981849611063
$data[9,0] = "VBL"
$data[9,1] = 490.0
$data[9,2] = 663.6
$data[9,3] = 73.83966244725738
$data[9,4] = 419.55
$data[9,5] = 116.7918007388869
$data[10,0] = "ABB"
$data[10,1] = 5167.5
$data[10,2] = 6948.4
$data[10,3] = 74.36963905359508
$data[10,4] = 4684.45
$data[10,5] = 110.31177619571135
$data[11,0] = "PFC"
$data[11,1] = 355.05
$data[11,2] = 474.85
$data[11,3] = 74.77098030957144
$data[11,4] = 329.9
$data[11,5] = 107.62352227947864
$data[12,0] = "DMART"
$data[12,1] = 3776.0
$data[12,2] = 4949.5
$data[12,3] = 76.29053439741388
$data[12,4] = 3340.0
$data[12,5] = 113.05389221556887
$data[13,0] = "BAJAJHLDNG"
$data[13,1] = 11305.0
$data[13,2] = 14763.0
$data[13,3] = 76.57657657657658
$data[13,4] = 10245.1
$data[13,5] = 110.3454334267113
$data[14,0] = "DLF"
$data[14,1] = 689.0
$data[14,2] = 886.8
$data[14,3] = 77.69508344609832
$data[14,4] = 601.2
$data[14,5] = 114.60412508316699
$data[15,0] = "ADANIPOWER"
$data[15,1] = 142.89
$data[15,2] = 182.7
$data[15,3] = 78.21018062397373
$data[15,4] = 89.0
$data[15,5] = 160.55056179775278
$data[16,0] = "IRFC"
$data[16,1] = 124.75
$data[16,2] = 158.0
$data[16,3] = 78.95569620253164
$data[16,4] = 108.04
$data[16,5] = 115.46649389115142
$data[17,0] = "HYUNDAI"
$data[17,1] = 2295.4
$data[17,2] = 2890.0
$data[17,3] = 79.42560553633218
$data[17,4] = 1541.7
$data[17,5] = 148.88759161964066
$data[18,0] = "CGPOWER"
$data[18,1] = 647.35
$data[18,2] = 797.55
$data[18,3] = 81.16732493260612
$data[18,4] = 517.7
$data[18,5] = 125.04346146416843
$data[19,0] = "SHREECEM"
$data[19,1] = 26635.0
$data[19,2] = 32490.0
$data[19,3] = 81.97907048322561
$data[19,4] = 24817.8
$data[19,5] = 107.32216393072714
$data[20,0] = "INDHOTEL"
$data[20,1] = 737.45
$data[20,2] = 894.9
$data[20,3] = 82.40585540283831
$data[20,4] = 672.6
$data[20,5] = 109.64168896818317
$data[21,0] = "HAVELLS"
$data[21,1] = 1421.4
$data[21,2] = 1721.2
$data[21,3] = 82.58191959098305
$data[21,4] = 1380.0
$data[21,5] = 103.0
$data[22,0] = "GAIL"
$data[22,1] = 172.0
$data[22,2] = 202.79
$data[22,3] = 84.81680556240447
$data[22,4] = 150.52
$data[22,5] = 114.27052883337761
$data[23,0] = "HAL"
$data[23,1] = 4386.0
$data[23,2] = 5165.0
$data[23,3] = 84.91771539206196
$data[23,4] = 3046.05
$data[23,5] = 143.9897572265721
$data[24,0] = "UNITDSPR"
$data[24,1] = 1444.8
$data[24,2] = 1700.0
$data[24,3] = 84.98823529411764
$data[24,4] = 1271.1
$data[24,5] = 113.6653292423885
$data[25,0] = "BOSCHLTD"
$data[25,1] = 36095.0
$data[25,2] = 41945.0
$data[25,3] = 86.0531648587436
$data[25,4] = 25921.6
$data[25,5] = 139.24680575273132
$data[26,0] = "ADANIGREEN"
$data[26,1] = 1015.0
$data[26,2] = 1177.55
$data[26,3] = 86.19591524776018
$data[26,4] = 758.0
$data[26,5] = 133.90501319261213
$data[27,0] = "ZYDUSLIFE"
$data[27,1] = 914.0
$data[27,2] = 1059.05
$data[27,3] = 86.30376280628866
$data[27,4] = 795.0
$data[27,5] = 114.96855345911949
$data[28,0] = "LICI"
$data[28,1] = 855.1
$data[28,2] = 980.0
$data[28,3] = 87.25510204081633
$data[28,4] = 715.3
$data[28,5] = 119.54424716902
$data[29,0] = "AMBUJACEM"
$data[29,1] = 556.5
$data[29,2] = 624.95
$data[29,3] = 89.04712376990159
$data[29,4] = 455.0
$data[29,5] = 122.3076923076923
$data[30,0] = "DIVISLAB"
$data[30,1] = 6398.0
$data[30,2] = 7071.5
$data[30,3] = 90.47585377925475
$data[30,4] = 4955.0
$data[30,5] = 129.1220988900101
$data[31,0] = "TATAPOWER"
$data[31,1] = 379.4
$data[31,2] = 416.8
$data[31,3] = 91.02687140115162
$data[31,4] = 326.35
$data[31,5] = 116.25555385322505
$data[32,0] = "HINDZINC"
$data[32,1] = 612.3
$data[32,2] = 656.35
$data[32,3] = 93.28864173078387
$data[32,4] = 378.15
$data[32,5] = 161.91987306624355
$data[33,0] = "GODREJCP"
$data[33,1] = 1225.0
$data[33,2] = 1309.0
$data[33,3] = 93.58288770053476
$data[33,4] = 979.5
$data[33,5] = 125.06380806533946
$data[34,0] = "PIDILITIND"
$data[34,1] = 1482.1
$data[34,2] = 1574.95
$data[34,3] = 94.1045747484047
$data[34,4] = 1311.1
$data[34,5] = 113.04248341087637
$data[35,0] = "ICICIGI"
$data[35,1] = 1961.0
$data[35,2] = 2068.7
$data[35,3] = 94.79383187509065
$data[35,4] = 1613.7
$data[35,5] = 121.52196814773501
$data[36,0] = "LTIM"
$data[36,1] = 6060.0
$data[36,2] = 6380.0
$data[36,3] = 94.98432601880879
$data[36,4] = 3802.0
$data[36,5] = 159.3897948448185
$data[37,0] = "BRITANNIA"
$data[37,1] = 6027.0
$data[37,2] = 6336.0
$data[37,3] = 95.12310606060606
$data[37,4] = 4506.0
$data[37,5] = 133.75499334221038
$data[38,0] = "IOC"
$data[38,1] = 166.34
$data[38,2] = 174.5
$data[38,3] = 95.32378223495702
$data[38,4] = 110.72
$data[38,5] = 150.23482658959537
$data[39,0] = "CHOLAFIN"
$data[39,1] = 1700.1
$data[39,2] = 1782.0
$data[39,3] = 95.4040404040404
$data[39,4] = 1169.8
$data[39,5] = 145.33253547614976
$data[40,0] = "JINDALSTEL"
$data[40,1] = 1050.0
$data[40,2] = 1098.0
$data[40,3] = 95.62841530054645
$data[40,4] = 723.35
$data[40,5] = 145.15794566945462
$data[41,0] = "PNB"
$data[41,1] = 123.7
$data[41,2] = 127.8
$data[41,3] = 96.79186228482003
$data[41,4] = 85.46
$data[41,5] = 144.74608003744444
$data[42,0] = "BANKBARODA"
$data[42,1] = 295.9
$data[42,2] = 303.95
$data[42,3] = 97.35153808192136
$data[42,4] = 190.7
$data[42,5] = 155.1651809124279
$data[43,0] = "MOTHERSON"
$data[43,1] = 119.99
$data[43,2] = 122.8
$data[43,3] = 97.71172638436482
$data[43,4] = 71.5
$data[43,5] = 167.8181818181818
$data[44,0] = "VEDL"
$data[44,1] = 602.4
$data[44,2] = 616.0
$data[44,3] = 97.79220779220779
$data[44,4] = 363.0
$data[44,5] = 165.9504132231405
$data[45,0] = "ADANIENSOL"
$data[45,1] = 1027.0
$data[45,2] = 1050.0
$data[45,3] = 97.80952380952381
$data[45,4] = 639.45
$data[45,5] = 160.60677144420984
$data[46,0] = "CANBK"
$data[46,1] = 154.8
$data[46,2] = 158.0
$data[46,3] = 97.9746835443038
$data[46,4] = 78.6
$data[46,5] = 196.94656488549623
$data[47,0] = "TORNTPHARM"
$data[47,1] = 3846.6
$data[47,2] = 3882.2
$data[47,3] = 99.08299417855856
$data[47,4] = 2886.45
$data[47,5] = 133.26404406797278
$data[48,0] = "BPCL"
$data[48,1] = 383.45
$data[48,2] = 386.0
$data[48,3] = 99.33937823834196
$data[48,4] = 234.01
$data[48,5] = 163.86051878124866
$data[49,0] = "TVSMOTOR"
$data[49,1] = 3720.0
$data[49,2] = 3734.9
$data[49,3] = 99.60106026935125
$data[49,4] = 2171.4
$data[49,5] = 171.31804365846918

$ws.Range("A2:F51").Value = $data

# Restore the originally active sheet/tab so this sheet addition is the only visible change
$wb.Worksheets.Item($origActiveSheetName).Activate()
